# KAMUS MSIC.xlsx — add a new dictionary entry row.
#
# The existing row 7 ("e-hailing") maps the term to MSIC code 49225
# ("Sewa kereta dengan pemandu") under Section H - "Pengangkutan Dan
# Penyimpanan". This change inserts a sibling row right below it for the
# Malay-language synonym "e-panggilan", carrying over the exact same
# MSIC code / description / notes / section columns (B-F), while column A
# gets the new term.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8 (pushes old row 8..105 down to 9..106).
# Excel's default Insert behaviour copies formatting from the row above,
# which is exactly what we want here (same style ids as row 7).
$ws.Rows.Item(8).Insert()

# Column A: the new term itself.
$ws.Cells.Item(8, 1).Value = "e-panggilan"

# Columns B-F: duplicate row 7's MSIC code / description / notes / section
# values verbatim (this reuses the existing shared strings instead of
# minting new ones).
$ws.Cells.Item(8, 2).Value = $ws.Cells.Item(7, 2).Text
$ws.Cells.Item(8, 3).Value = $ws.Cells.Item(7, 3).Text
$ws.Cells.Item(8, 4).Value = $ws.Cells.Item(7, 4).Text
$ws.Cells.Item(8, 5).Value = $ws.Cells.Item(7, 5).Text
$ws.Cells.Item(8, 6).Value = $ws.Cells.Item(7, 6).Text

# Match row 7's (taller, wrapped-note) row height instead of the
# auto-fit height Excel would otherwise pick for the new row.
$ws.Rows.Item(8).RowHeight = $ws.Rows.Item(7).RowHeight

# Leave the cursor where the author's last save shows it.
$ws.Range("C58").Select()
